$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading_percent values for rows 2-25 (columns B,D,E,F,G,I,J,K,L,N)
$newValues = @{
    2 = @{ "B"=16.02215377896317; "D"=3.960574696401489; "E"=8.962828263891319; "F"=49.9221668532283; "G"=3.764367853488253; "I"=33.36126275062082; "J"=9.007884358450323; "K"=18.83968362875145; "L"=13.75599988110503; "N"=23.78713699837153 }
    3 = @{ "B"=15.94831283836602; "D"=3.938544846372974; "E"=8.940634362825559; "F"=49.90977056915472; "G"=3.767743117104273; "I"=33.42382773433575; "J"=9.005588623716003; "K"=18.68483624650251; "L"=13.74084287428061; "N"=23.84462511146799 }
    4 = @{ "B"=15.90625816759909; "D"=3.924672697802881; "E"=8.926688415411974; "F"=49.91230821461373; "G"=3.769924619531288; "I"=33.46704985798072; "J"=9.004209262063267; "K"=18.59400664939239; "L"=13.73418770973137; "N"=23.88186005853999 }
    5 = @{ "B"=15.88996043417543; "D"=3.918932649535146; "E"=8.920924845494586; "F"=49.91589418597442; "G"=3.770841125319543; "I"=33.48587069359929; "J"=9.003654674384707; "K"=18.55809688647252; "L"=13.73214442155381; "N"=23.89752141518379 }
    6 = @{ "B"=15.88730532555008; "D"=3.917974259188505; "E"=8.919962938727638; "F"=49.91664372779439; "G"=3.770994975767217; "I"=33.48906876792431; "J"=9.003563036759338; "K"=18.55220178864648; "L"=13.73184557728028; "N"=23.9001514560998 }
    7 = @{ "B"=15.90603495220403; "D"=3.924595638011091; "E"=8.926611012156391; "F"=49.91234624528984; "G"=3.769936868269156; "I"=33.46729879519059; "J"=9.004201752383182; "K"=18.59351784164237; "L"=13.7341574431351; "N"=23.88206929699705 }
    8 = @{ "B"=15.99602062759356; "D"=3.953050499654995; "E"=8.955241780701691; "F"=49.91578645283391; "G"=3.765509064290702; "I"=33.38183672283581; "J"=9.007086461435881; "K"=18.78543296775857; "L"=13.75022482453754; "N"=23.80655731570936 }
    9 = @{ "B"=16.19786995537265; "D"=4.006112988537727; "E"=9.008872788648342; "F"=50.00301923336342; "G"=3.7576871836854; "I"=33.25244287349561; "J"=9.012987497898164; "K"=19.19379858420195; "L"=13.80266262658217; "N"=23.67381387469858 }
    10 = @{ "B"=16.36068244613568; "D"=4.043439362346437; "E"=9.046766756148061; "F"=50.11604686446381; "G"=3.752459136163049; "I"=33.1807372519142; "J"=9.017478198353592; "K"=19.51097491453263; "L"=13.85377980071124; "N"=23.58558791647628 }
    11 = @{ "B"=16.43767746726408; "D"=4.060062142409282; "E"=9.063683210371044; "F"=50.17803592002584; "G"=3.750192064617766; "I"=33.1532029897752; "J"=9.019556243840693; "K"=19.65845820296797; "L"=13.87972360468021; "N"=23.5474611989704 }
    12 = @{ "B"=16.46723482387681; "D"=4.066305500123091; "E"=9.07004311165421; "F"=50.20302240210576; "G"=3.749349471718253; "I"=33.14350864384018; "J"=9.0203483127959; "K"=19.71472015123889; "L"=13.88993017794492; "N"=23.53331163171149 }
    13 = @{ "B"=16.4608516155071; "D"=4.064963163686399; "E"=9.068675434461799; "F"=50.19757398229216; "G"=3.749530233564832; "I"=33.14556391054338; "J"=9.020177495627427; "K"=19.70258550379538; "L"=13.8877150891005; "N"=23.53634618465538 }
    14 = @{ "B"=16.44010125128717; "D"=4.0605768157221; "E"=9.064207363523865; "F"=50.18006130808755; "G"=3.750122425898597; "I"=33.15239074840684; "J"=9.019621303898324; "K"=19.66307889467141; "L"=13.88055567354239; "N"=23.54629133260329 }
    15 = @{ "B"=16.42744265992489; "D"=4.057883367409589; "E"=9.06146455930466; "F"=50.16953101734109; "G"=3.75048722846693; "I"=33.15666777800355; "J"=9.019281295021807; "K"=19.63893237664659; "L"=13.87621995553203; "N"=23.55242053926128 }
    16 = @{ "B"=16.35570788467276; "D"=4.042345840196975; "E"=9.045654775602092; "F"=50.11220784938597; "G"=3.752609524275441; "I"=33.18263909751992; "J"=9.017343109769467; "K"=19.5013969136884; "L"=13.8521380609471; "N"=23.58811994306684 }
    17 = @{ "B"=16.31243725571867; "D"=4.032722613997093; "E"=9.035873781801959; "F"=50.07974473087232; "G"=3.753939897989678; "I"=33.19987484417484; "J"=9.016163180437596; "K"=19.41780947478147; "L"=13.83805048781688; "N"=23.61053422231998 }
    18 = @{ "B"=16.28782626499258; "D"=4.02715396809825; "E"=9.030217717101182; "F"=50.06206846601914; "G"=3.754715565266151; "I"=33.21026695640875; "J"=9.015487824014276; "K"=19.37003605676667; "L"=13.83020122223744; "N"=23.62361532782434 }
    19 = @{ "B"=16.27954161506776; "D"=4.025262744849878; "E"=9.028297459235805; "F"=50.05625480276118; "G"=3.754979994167617; "I"=33.2138677068001; "J"=9.015259725126288; "K"=19.35391440019506; "L"=13.82758727512585; "N"=23.62807684312114 }
    20 = @{ "B"=16.31701495198326; "D"=4.0337505021707; "E"=9.036918118776477; "F"=50.08309748940229; "G"=3.753797194344171; "I"=33.19799053175694; "J"=9.016288442882887; "K"=19.42667638780931; "L"=13.83952392662003; "N"=23.60812862662755 }
    21 = @{ "B"=16.44618541901568; "D"=4.061866585235301; "E"=9.065520990390654; "F"=50.18516422289429; "G"=3.749948053917478; "I"=33.15036565881832; "J"=9.019784530130408; "K"=19.67467210289913; "L"=13.88264823332815; "N"=23.54336238407994 }
    22 = @{ "B"=16.53293327018288; "D"=4.079942903573249; "E"=9.083946660120468; "F"=50.26068255415499; "G"=3.747525041091063; "I"=33.12350875814152; "J"=9.022099536416176; "K"=19.83913831325811; "L"=13.91305786358837; "N"=23.50271348807834 }
    23 = @{ "B"=16.48642826224832; "D"=4.070322595913854; "E"=9.074136964645103; "F"=50.21957364034333; "G"=3.748809803756476; "I"=33.13745188310989; "J"=9.02086118691199; "K"=19.75115675151451; "L"=13.89662569269891; "N"=23.52425505057146 }
    24 = @{ "B"=16.31494454580464; "D"=4.033285906478594; "E"=9.036446076372929; "F"=50.08157863159608; "G"=3.75386167695937; "I"=33.1988409247978; "J"=9.016231802400926; "K"=19.42266677412824; "L"=13.83885700627936; "N"=23.60921558947363 }
    25 = @{ "B"=16.14064209716727; "D"=3.992050162721315; "E"=8.994630823440172; "F"=49.9708114386872; "G"=3.759711675315176; "I"=33.28335002245119; "J"=9.011364317879108; "K"=19.08013571988089; "L"=13.78625175135595; "N"=23.70808779723474 }
}

foreach ($row in $newValues.Keys) {
    $rowData = $newValues[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
